$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlineShapeImage($inlineShape, $newName) {
    # InlineShape has no writable "Name"/filename property of its own in the
    # Word object model (only Shape does), so round-trip the picture through
    # a floating Shape - which does expose a settable .Name - and convert it
    # back to an inline picture afterwards so the layout/wrapping is
    # unchanged.
    $floating = $inlineShape.ConvertToShape()
    $floating.Name = $newName
    $floating.ConvertToInlineShape() | Out-Null
}

# Footer (default / "Primary") - Pearson Edexcel logo: image1.png -> image2.png
$footerPrimary = $sec.Footers(1)
Rename-InlineShapeImage $footerPrimary.Range.InlineShapes(1) "image2.png"

# Footer (first page) - Pearson Edexcel logo: image1.png -> image2.png
$footerFirst = $sec.Footers(2)
Rename-InlineShapeImage $footerFirst.Range.InlineShapes(1) "image2.png"

# Header (first page) - BTEC logo: image2.jpg -> image1.jpg
$headerFirst = $sec.Headers(2)
Rename-InlineShapeImage $headerFirst.Range.InlineShapes(1) "image1.jpg"
